$d = $word.ActiveDocument

# Locate the run "...)}" that closes the {m: ... } field text.
$rng = $d.Content
$found = $rng.Find.Execute(")}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the ')}' text to split"
}

$match = $rng.Duplicate
$matchEnd = $match.End

# Step 1: split the single run containing ")}" into two runs - one for ")"
# and one for "}" - by toggling (and reverting) a character formatting
# property on just the last character. Reverting the value keeps the
# resulting run properties identical to the original, while still forcing
# Word to materialize a separate run boundary right before "}".
$closeBrace = $d.Range($matchEnd - 1, $matchEnd)
$closeBrace.Font.Bold = 1
$closeBrace.Font.Bold = 0

# Step 2: rewrite just that new "}" run so it matches the expected output
# of the field rewriter: a fresh run (no rsid carried over) whose text is
# marked with xml:space="preserve".
$closeBrace2 = $d.Range($matchEnd - 1, $matchEnd)
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$closeBrace2.InsertXML($xmlFrag)
